$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: replace CA901 (Ebony Polish) duplicate with CA701 Premium Satin Black data ---
$ws.Range("A6").Value = "https://store.kawaius.com/products/productdetail/CA701+Premium+Satin+Black+Digital+Piano/part_number=C11-CA701SB/1772.0.1.1.59003.0.0.0.0?pp=8&"
$ws.Range("B6").Value = "CA701 Premium Satin Black Digital Piano"
$ws.Range("C6").Value = 6099
$ws.Range("D6").Value = 4899
$ws.Range("A6:D6").Style = "Normal"

# --- Row 7 (new): CA701 Premium Satin White ---
$ws.Range("A7").Value = "https://store.kawaius.com/products/productdetail/CA701+Premium+Satin+White+Digital+Piano/part_number=C11-CA701WH/1772.0.1.1.59003.0.0.0.0?s=part_number&part_number_d=ASC&part_number_c=part_number&t=1&i=all&"
$ws.Range("B7").Value = "CA701 Premium Satin White Digital Piano"
$ws.Range("C7").Value = 6099
$ws.Range("D7").Value = 4899
$ws.Range("A7:D7").Style = "Normal"
$ws.Rows.Item(7).RowHeight = 14.25

# --- Row 8: replace short CN301SB url with mp7se stage piano url; zero-out MSRP/MAP ---
$ws.Range("A8").Value = "https://store.kawaius.com/p/mp7se-professional-stage-piano?pp=8"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("A8").Style = "Normal"
$ws.Range("C8:D8").Style = "Normal"

# --- Row 9: replace CA901 duplicate url with es110 pedal bar url; zero-out MSRP/MAP ---
$ws.Range("A9").Value = "https://store.kawaius.com/p/es110-black-triple-pedal-bar?pp=8"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("A9").Style = "Normal"
$ws.Range("C9:D9").Style = "Normal"

# --- Row 10: replace gl-10 duplicate url with k-200 upright piano ---
$ws.Range("A10").Value = "https://kawaius.com/product/k-200/"
$ws.Range("B10").Value = "Kawai K-200 Upright Piano"
$ws.Range("C10").Value = 8095
$ws.Range("D10").Value = 0
$ws.Range("A10:D10").Style = "Normal"

# --- Row 11: replace cn301 duplicate url with gl-40 grand piano ---
$ws.Range("A11").Value = "https://kawaius.com/product/gl-40/"
$ws.Range("B11").Value = "Kawai GL-40 Grand Piano"
$ws.Range("C11").Value = 38895
$ws.Range("D11").Value = 0
$ws.Range("A11:D11").Style = "Normal"

# --- Row 12 (new): k-500 upright piano ---
$ws.Range("A12").Value = "https://kawaius.com/product/k-500/"
$ws.Range("B12").Value = "Kawai K-500 Upright Piano"
$ws.Range("C12").Value = 16095
$ws.Range("D12").Value = 0
$ws.Range("A12:D12").Style = "Normal"
$ws.Rows.Item(12).RowHeight = 14.25

# --- Leftover formatted-but-empty rows from the sheet clear pass ---
$ws.Rows.Item(14).RowHeight = 14.25
$ws.Rows.Item(18).RowHeight = 14.25
$ws.Rows.Item(19).RowHeight = 14.25
$ws.Rows.Item(21).RowHeight = 14.25

# --- Selection moves to A9 ---
$null = $ws.Range("A9").Select()
